$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147; existing rows 147-202 shift down to 148-203.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new weekly data point.
$ws.Range("A147").Value = 7
$ws.Range("B147").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C147").Value = "Ñuble"
$ws.Range("D147").Value = 45009
$ws.Range("E147").Value = 16
$ws.Range("F147").Value = 100112040
$ws.Range("G147").Value = "Cilantro"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 530
$ws.Range("K147").Value = 1000
$ws.Range("L147").Value = 1200
$ws.Range("M147").Value = 1106
$ws.Range("N147").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O147").Value = "Provincia de Diguillín"
$ws.Range("P147").Value = 1106
$ws.Range("Q147").Value = 1
$ws.Range("R147").Value = "Hortaliza"
